$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.895.99"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "3.839.15"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "637.13"
$ws.Range("E5").Value = "  +6.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.78"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "3.836.36"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.68"
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "4.480.50"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "3.984.23"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "69.803.08"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.79"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.85"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  +3.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "3.986.20"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.72"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.34"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.36"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "3.781.60"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.08"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  +3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.151"
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.48"
$ws.Range("E40").Value = "  +6.82%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.984"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.58"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.05"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.48"
$ws.Range("E51").Value = "  +1.56%  "
